{"js": "// Replace the date line and each \"A\u00f7B=\" problem text with its updated\n// value, in document order. A plain global find/replace is unsafe here\n// because one of the new values (\"88\u00f72=\") duplicates an old value that\n// appears elsewhere in the document, so we walk paragraphs positionally\n// instead and only touch the ones that actually hold one of the known\n// old values.\nconst replacements = [\n  \"2023-12-02 Saturday\",\n  \"79\u00f78=\",\n  \"73\u00f73=\",\n  \"21\u00f77=\",\n  \"33\u00f72=\",\n  \"26\u00f76=\",\n  \"45\u00f76=\",\n  \"37\u00f78=\",\n  \"67\u00f73=\",\n  \"41\u00f78=\",\n  \"36\u00f75=\",\n  \"99\u00f78=\",\n  \"61\u00f72=\",\n  \"78\u00f73=\",\n  \"23\u00f79=\",\n  \"88\u00f72=\",\n  \"65\u00f78=\",\n  \"68\u00f78=\",\n  \"73\u00f77=\",\n  \"96\u00f72=\",\n  \"80\u00f77=\",\n  \"99\u00f79=\",\n  \"51\u00f79=\",\n  \"20\u00f76=\",\n  \"78\u00f78=\",\n  \"70\u00f72=\",\n];\n\nconst oldValues = new Set([\n  \"2023-12-01 Friday\",\n  \"42\u00f74=\",\n  \"56\u00f75=\",\n  \"14\u00f73=\",\n  \"45\u00f77=\",\n  \"88\u00f75=\",\n  \"88\u00f72=\",\n  \"48\u00f72=\",\n  \"38\u00f75=\",\n  \"62\u00f77=\",\n  \"81\u00f78=\",\n  \"10\u00f79=\",\n  \"82\u00f74=\",\n  \"77\u00f79=\",\n  \"78\u00f77=\",\n  \"19\u00f79=\",\n  \"22\u00f72=\",\n  \"39\u00f76=\",\n  \"35\u00f74=\",\n  \"66\u00f75=\",\n  \"55\u00f73=\",\n  \"17\u00f78=\",\n  \"29\u00f76=\",\n  \"72\u00f78=\",\n  \"21\u00f78=\",\n  \"50\u00f79=\",\n]);\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet idx = 0;\nfor (const para of paragraphs.items) {\n  if (idx >= replacements.length) break;\n  if (oldValues.has(para.text)) {\n    para.insertText(replacements[idx], \"Replace\");\n    idx++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the date line and each \"A\u00f7B=\" problem text with its updated\n# value, in document order. A blind ExecuteFind/Replace-All is unsafe\n# here because one of the new values (\"88\u00f72=\") duplicates an old value\n# that appears elsewhere in the document, so instead we walk the\n# paragraphs positionally and only touch the ones that still hold one\n# of the known old values, consuming the replacement list in order.\n\n$d = $word.ActiveDocument\n\n$oldValues = @(\n    \"2023-12-01 Friday\",\n    \"42\u00f74=\",\n    \"56\u00f75=\",\n    \"14\u00f73=\",\n    \"45\u00f77=\",\n    \"88\u00f75=\",\n    \"88\u00f72=\",\n    \"48\u00f72=\",\n    \"38\u00f75=\",\n    \"62\u00f77=\",\n    \"81\u00f78=\",\n    \"10\u00f79=\",\n    \"82\u00f74=\",\n    \"77\u00f79=\",\n    \"78\u00f77=\",\n    \"19\u00f79=\",\n    \"22\u00f72=\",\n    \"39\u00f76=\",\n    \"35\u00f74=\",\n    \"66\u00f75=\",\n    \"55\u00f73=\",\n    \"17\u00f78=\",\n    \"29\u00f76=\",\n    \"72\u00f78=\",\n    \"21\u00f78=\",\n    \"50\u00f79=\"\n)\n\n$newValues = @(\n    \"2023-12-02 Saturday\",\n    \"79\u00f78=\",\n    \"73\u00f73=\",\n    \"21\u00f77=\",\n    \"33\u00f72=\",\n    \"26\u00f76=\",\n    \"45\u00f76=\",\n    \"37\u00f78=\",\n    \"67\u00f73=\",\n    \"41\u00f78=\",\n    \"36\u00f75=\",\n    \"99\u00f78=\",\n    \"61\u00f72=\",\n    \"78\u00f73=\",\n    \"23\u00f79=\",\n    \"88\u00f72=\",\n    \"65\u00f78=\",\n    \"68\u00f78=\",\n    \"73\u00f77=\",\n    \"96\u00f72=\",\n    \"80\u00f77=\",\n    \"99\u00f79=\",\n    \"51\u00f79=\",\n    \"20\u00f76=\",\n    \"78\u00f78=\",\n    \"70\u00f72=\"\n)\n\n$idx = 0\nforeach ($p in $d.Paragraphs) {\n    if ($idx -ge $newValues.Length) { break }\n    $r = $p.Range\n    $text = $r.Text.TrimEnd([char]13, [char]7)\n    if ($oldValues -contains $text) {\n        $r.Text = $newValues[$idx]\n        $idx++\n    }\n}\n"}
